# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) on several Leve-profit sheets with freshly
# fetched values. Mirrors an automated scraper run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 737.5
$ws.Range("J12").Value = 766.6667
$ws.Range("L12").Value = 766.6667
$ws.Range("N12").Value = -1106.6667

# Row 51
$ws.Range("H51").Value = 2963.5454
$ws.Range("I51").Value = 2799
$ws.Range("J51").Value = 2980
$ws.Range("K51").Value = 2799
$ws.Range("L51").Value = 2980
$ws.Range("M51").Value = -2315
$ws.Range("N51").Value = -3948

# Row 95
$ws.Range("H95").Value = 45000
$ws.Range("J95").Value = 45000
$ws.Range("L95").Value = 45000
$ws.Range("N95").Value = -50492

# Row 135
$ws.Range("H135").Value = 912.1667
$ws.Range("I135").Value = 732.7692
$ws.Range("J135").Value = 1378.6
$ws.Range("K135").Value = 6594.922799999999
$ws.Range("L135").Value = 12407.4
$ws.Range("M135").Value = -4059.922799999999
$ws.Range("N135").Value = -17477.4

# Row 138
$ws.Range("H138").Value = 2452.0952
$ws.Range("I138").Value = 856.9167
$ws.Range("J138").Value = 4579
$ws.Range("K138").Value = 2570.7501
$ws.Range("L138").Value = 13737
$ws.Range("M138").Value = 2569.2499
$ws.Range("N138").Value = -24017

$ws = $wb.Worksheets.Item("ARM")
# Row 46
$ws.Range("H46").Value = 3500
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Row 61
$ws.Range("H61").Value = 1680.7778
$ws.Range("I61").Value = 1534.625
$ws.Range("K61").Value = 1534.625
$ws.Range("M61").Value = -1322.625

# Row 96
$ws.Range("H96").Value = 2025972.6
$ws.Range("J96").Value = 2025972.6
$ws.Range("L96").Value = 2025972.6
$ws.Range("N96").Value = -2031464.6

# Row 132
$ws.Range("H132").Value = 3246.4546
$ws.Range("I132").Value = 3344.375
$ws.Range("J132").Value = 2985.3333
$ws.Range("K132").Value = 10033.125
$ws.Range("L132").Value = 8955.999899999999
$ws.Range("M132").Value = -7503.125
$ws.Range("N132").Value = -14015.9999

# Row 133
$ws.Range("H133").Value = 230125
$ws.Range("J133").Value = 230125
$ws.Range("L133").Value = 230125
$ws.Range("N133").Value = -235185

# Row 136
$ws.Range("H136").Value = 1680.7778
$ws.Range("I136").Value = 1534.625
$ws.Range("K136").Value = 4603.875
$ws.Range("M136").Value = -2053.875

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5841.643
$ws.Range("I134").Value = 1706.3846
$ws.Range("K134").Value = 5119.1538
$ws.Range("M134").Value = -2584.1538

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 551.8
$ws.Range("I16").Value = 640
$ws.Range("K16").Value = 640
$ws.Range("M16").Value = -353

# Row 48
$ws.Range("H48").Value = 49999
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -50951

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 58
$ws.Range("H58").Value = 6372.5
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

# Row 113
$ws.Range("H113").Value = 551.8
$ws.Range("I113").Value = 640
$ws.Range("K113").Value = 640
$ws.Range("M113").Value = 1530

# Row 132
$ws.Range("H132").Value = 2127.0454
$ws.Range("I132").Value = 1863.6111
$ws.Range("K132").Value = 5590.8333
$ws.Range("M132").Value = -3060.8333

# Row 134
$ws.Range("H134").Value = 2148.55
$ws.Range("I134").Value = 1345.6111
$ws.Range("K134").Value = 4036.8333
$ws.Range("M134").Value = -1501.8333

# Row 136
$ws.Range("H136").Value = 6372.5
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1391279.5
$ws.Range("I4").Value = 1391279.5
$ws.Range("K4").Value = 4173838.5
$ws.Range("M4").Value = -4173726.5

# Row 12
$ws.Range("H12").Value = 191.64285
$ws.Range("I12").Value = 39
$ws.Range("J12").Value = 203.38461
$ws.Range("K12").Value = 117
$ws.Range("L12").Value = 610.15383
$ws.Range("M12").Value = 56
$ws.Range("N12").Value = -956.15383

# Row 23
$ws.Range("H23").Value = 347.22223
$ws.Range("I23").Value = 347.5
$ws.Range("J23").Value = 347
$ws.Range("K23").Value = 1042.5
$ws.Range("L23").Value = 1041
$ws.Range("M23").Value = -807.5
$ws.Range("N23").Value = -1511

# Row 103
$ws.Range("H103").Value = 1614.4286
$ws.Range("J103").Value = 1724.8462
$ws.Range("L103").Value = 5174.5386
$ws.Range("N103").Value = -6932.5386

# Row 122
$ws.Range("H122").Value = 611.7619
$ws.Range("I122").Value = 374.4
$ws.Range("K122").Value = 3369.6
$ws.Range("M122").Value = -919.5999999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Range("H98").Value = 8603.909
$ws.Range("J98").Value = 8564.299999999999
$ws.Range("L98").Value = 8564.299999999999
$ws.Range("N98").Value = -14554.3

# Row 132
$ws.Range("H132").Value = 117210.89
$ws.Range("I132").Value = 172616.33
$ws.Range("K132").Value = 517848.99
$ws.Range("M132").Value = -515318.99

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 4250
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 4444.4443
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 13333.3329
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -18433.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7900
$ws.Range("J62").Value = 8800
$ws.Range("L62").Value = 8800
$ws.Range("N62").Value = -10048

# Row 65
$ws.Range("H65").Value = 7900
$ws.Range("J65").Value = 8800
$ws.Range("L65").Value = 44000
$ws.Range("N65").Value = -50240

# Row 96
$ws.Range("H96").Value = 2079
$ws.Range("I96").Value = 2131.6667
$ws.Range("K96").Value = 2131.6667
$ws.Range("M96").Value = -758.6667000000002

# Row 113
$ws.Range("H113").Value = 1204.5
$ws.Range("I113").Value = 1072.5
$ws.Range("J113").Value = 1270.5
$ws.Range("K113").Value = 3217.5
$ws.Range("L113").Value = 3811.5
$ws.Range("M113").Value = -1047.5
$ws.Range("N113").Value = -8151.5

# Row 136
$ws.Range("H136").Value = 2362.634
$ws.Range("I136").Value = 1711.5
$ws.Range("J136").Value = 4677.778
$ws.Range("K136").Value = 5134.5
$ws.Range("L136").Value = 14033.334
$ws.Range("M136").Value = -2584.5
